$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Activate()
$ws.Rows.Item(6).Resize(2).Insert()
$ws.Range("A6").Value = "Prior distribution for fluxes (uniform or normal)"
$ws.Range("B6").Value = "normal"
$ws.Range("A7").Value = "Prior distribution for thermodynamic quantities (uniform or normal)"
$ws.Range("B7").Value = "normal"
$ws.Rows.Item(1048577).Resize(2).Delete()
$ws.Range("A6:B7").Select()
